# Updates Coin/Link/Price/Volume(1h) columns with the latest cryptos snapshot.
# Numeric-looking "Price" text values (e.g. "1.004") are forced back to plain
# text storage (matching the source data, which stores prices as strings) by
# toggling the number format to Text for the assignment, then restoring the
# cell style to Normal so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.557.94"
$ws.Range("E2").Value = "  +5.47%  "

$ws.Range("D3").Value = "1.723.53"
$ws.Range("E3").Value = "  +4.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5390"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2681"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.96%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06617"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.13%  "

$ws.Range("E10").Value = "  +6.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07716"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.616"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.962.25"
$ws.Range("E13").Value = "  +4.24%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.658.10"
$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5873"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.59%  "

$ws.Range("D16").Value = "0.0₅8319"
$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("E17").Value = "  +3.88%  "

$ws.Range("D18").Value = "27.577.86"
$ws.Range("E18").Value = "  +5.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "221.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +15.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.740"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.89%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.101"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.47%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.27"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.696"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.14%  "

$ws.Range("E27").Value = "  +3.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.406"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("E29").Value = "  +4.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05564"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.304"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.552"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.458"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.661"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9637"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.821"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.48%  "

$ws.Range("E37").Value = "  +1.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5959"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01649"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.31%  "

$ws.Range("E40").Value = "  +1.09%  "

$ws.Range("D41").Value = "1.056.50"
$ws.Range("E41").Value = "  +2.64%  "

$ws.Range("E42").Value = "  +2.22%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "

$ws.Range("D45").Value = "1.867.87"
$ws.Range("E45").Value = "  +4.14%  "

$ws.Range("E46").Value = "  +12.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.11"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.191"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "

$ws.Range("E49").Value = "  +2.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.006"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05277"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "
